# Auto-update draw results: append the 2025-12-09 Pick 3 draw as a new
# row (84) at the bottom of the results table on the active sheet.
#
# The source data is authored as plain text for every column (even the
# date-looking and all-digit-looking ones), so each value is entered with
# a leading apostrophe to force Excel's normal "smart" type detection
# (dates / numbers) to back off and keep the cell as text - matching the
# existing rows, which are all text cells too. The style is then reset to
# "Normal" so the new cells don't pick up a distinct quote-prefix format
# and stay visually identical to the rest of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Next empty row right after the current last data row (83 -> 84).
$row = $ws.UsedRange.Rows.Count + 1

$ws.Range("A$row").Value = "'2025-12-09"
$ws.Range("B$row").Value = "Pick 3"
$ws.Range("C$row").Value = "'251209"
$ws.Range("D$row").Value = "2-0-7"
$ws.Range("E$row").Value = "2025-12-09T21:42:45.126+04:00"

# Drop the quote-prefix formatting picked up from the apostrophes above so
# the new row keeps the same (default/"Normal") style as every other row.
$ws.Range("A$($row):E$row").Style = "Normal"

# Keep Excel's "number stored as text" warning suppressed for the whole
# table, now that it covers one more row.
$ws.Range("A1:E$row").Errors(9).Ignore = $true
